$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '33.877.32'
    'E2' = '  +9.20%  '
    'D3' = '1.782.07'
    'E3' = '  +5.68%  '
    'E4' = '  +0.06%  '
    'D5' = '224.97'
    'E5' = '  +2.02%  '
    'E6' = '  +4.30%  '
    'E7' = '  +0.08%  '
    'D8' = '30.77'
    'E8' = '  +4.88%  '
    'D9' = '46.52'
    'E9' = '  +3.89%  '
    'E10' = '  +4.19%  '
    'D11' = '0.0660'
    'E11' = '  +3.46%  '
    'E12' = '  +1.65%  '
    'D13' = '2.042.70'
    'E13' = '  +5.96%  '
    'D14' = '1.788.31'
    'E14' = '  +6.07%  '
    'D15' = '0.628'
    'E15' = '  +3.42%  '
    'D16' = '33.881.03'
    'E16' = '  +9.17%  '
    'D17' = '9.97'
    'E17' = '  -2.27%  '
    'D18' = '4.18'
    'E18' = '  +0.91%  '
    'D19' = '68.40'
    'E19' = '  +2.24%  '
    'D20' = '251.41'
    'E20' = '  +1.41%  '
    'E21' = '  +2.50%  '
    'D22' = '1.00'
    'E22' = '  +0.15%  '
    'D23' = '10.30'
    'E23' = '  +2.59%  '
    'D24' = '4.21'
    'E24' = '  -1.88%  '
    'E25' = '  +0.11%  '
    'D26' = '159.05'
    'E26' = '  +0.17%  '
    'D27' = '16.48'
    'E27' = '  +3.49%  '
    'D28' = '0.114'
    'E28' = '  +1.54%  '
    'E29' = '  +3.25%  '
    'D31' = '3.79'
    'E31' = '  +7.77%  '
    'E32' = '  +2.86%  '
    'E33' = '  +3.52%  '
    'E34' = '  +6.28%  '
    'D35' = '1.486.26'
    'E35' = '  -1.77%  '
    'E36' = '  +3.36%  '
    'E37' = '  +2.75%  '
    'E38' = '  +2.75%  '
    'D39' = '83.43'
    'E39' = '  +0.18%  '
    'E40' = '  +2.52%  '
    'E41' = '  +2.59%  '
    'E42' = '  +0.69%  '
    'D43' = '0.887'
    'E43' = '  +5.27%  '
    'D44' = '2.08'
    'E44' = '  +2.18%  '
    'D45' = '0.0508'
    'E45' = '  +0.90%  '
    'E46' = '  +2.90%  '
    'D47' = '1.939.62'
    'E47' = '  +6.40%  '
    'E48' = '  +3.57%  '
    'E49' = '  +0.10%  '
    'D50' = '11.90'
    'E50' = '  +15.76%  '
    'D51' = '50.74'
}

foreach ($addr in $changes.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $changes[$addr]
    $rng.Style = "Normal"
}
